# Applies the "Logging_1441" Lerum overview refresh:
#   - Every data row's "Forandrad" (C column) date bumps 46064 -> 46065
#   - The whole data block (rows 2..132) gets re-sorted into a new row
#     order by the upstream feed; a row's full content (A..Z) always
#     travels together with its case ("Beteckning"), it just lands on a
#     different row number. The mapping below is newRow -> oldRow.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowMap = @{ 2=2; 3=3; 4=5; 5=4; 6=6; 7=10; 8=8; 9=7; 10=9; 11=11; 12=12; 13=13; 14=14; 15=15; 16=16; 17=17; 18=18; 19=19; 20=20; 21=21; 22=22; 23=28; 24=26; 25=29; 26=30; 27=27; 28=23; 29=24; 30=31; 31=32; 32=62; 33=121; 34=71; 35=99; 36=25; 37=70; 38=96; 39=36; 40=118; 41=91; 42=92; 43=95; 44=126; 45=127; 46=129; 47=93; 48=41; 49=128; 50=48; 51=56; 52=59; 53=61; 54=60; 55=65; 56=106; 57=124; 58=67; 59=107; 60=131; 61=78; 62=38; 63=130; 64=82; 65=80; 66=46; 67=73; 68=68; 69=49; 70=47; 71=39; 72=69; 73=42; 74=120; 75=57; 76=53; 77=101; 78=94; 79=40; 80=35; 81=108; 82=33; 83=37; 84=116; 85=125; 86=132; 87=34; 88=43; 89=44; 90=45; 91=51; 92=52; 93=54; 94=58; 95=50; 96=55; 97=63; 98=64; 99=66; 100=72; 101=86; 102=79; 103=89; 104=85; 105=74; 106=84; 107=81; 108=75; 109=90; 110=83; 111=76; 112=77; 113=87; 114=88; 115=97; 116=98; 117=100; 118=102; 119=103; 120=105; 121=104; 122=111; 123=110; 124=109; 125=117; 126=114; 127=115; 128=112; 129=113; 130=119; 131=122; 132=123 }

$firstRow = 2
$lastRow = 132
$lastCol = "Z"

# 1) Snapshot every source row's formulas (A..Z) before any writes happen,
#    since several rows are both a source and a destination at different
#    row numbers.
$snapshots = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapshots[$r] = $ws.Range("A" + $r + ":" + $lastCol + $r).Formula
}

# 2) Write each destination row from its mapped source snapshot.
for ($newRow = $firstRow; $newRow -le $lastRow; $newRow++) {
    $oldRow = $rowMap[$newRow]
    $ws.Range("A" + $newRow + ":" + $lastCol + $newRow).Formula = $snapshots[$oldRow]
}

# 3) Bump the "Forandrad" column (C) from 46064 to 46065 for every data row.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Range("C" + $r).Value = 46065
}

# 4) The bulk Formula writes can disturb row heights (auto-fit kicking in
#    for the wrapped species-list column); restore the original fixed
#    15pt row height used throughout the data block.
for ($r = $firstRow; $r -le ($lastRow - 1); $r++) {
    $ws.Rows($r).RowHeight = 15
}
